$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 15
$ws.Range("I2").Value = 36
$ws.Range("J2").Value = 160
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 40
$ws.Range("N2").Value = 25
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 19
$ws.Range("T2").Value = 25
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 238
$ws.Range("X2").Value = 235
$ws.Range("Z2").Value = 7
$ws.Range("AA2").Value = 2
